# Reorganization and editing readme files
#
# 1) Slide 1 ("Seed Lab"): fill in the subtitle with the team member names.
# 2) Append 7 new "Title and Content" slides (layout index 2) with the
#    project outline content, in final presentation order.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - Subtitle: team member names
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange
$tr.Text = "Jonah "
[void]$tr.InsertAfter("Bertolino")
[void]$tr.InsertAfter(", Hunter Burnham, Joseph Kirby, Joel ")
[void]$tr.InsertAfter("Shorey")
[void]$tr.InsertAfter(",  Caden ")
[void]$tr.InsertAfter("Nubel")
$full = $subtitle.TextFrame.TextRange.Characters(1, $tr.Length)
$full.Font.Size = 20

# ---------------------------------------------------------------------
# Slide 2 - Project Overview
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Project Overview"

$tr = $s2.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "General Overview"
[void]$tr.InsertAfter("`rMeasures of success")
[void]$tr.InsertAfter("`rWell defined technical information (i.e. define ")
[void]$tr.InsertAfter("aruco")
[void]$tr.InsertAfter(" markers, PID or whatever ")
[void]$tr.InsertAfter("controllors")
[void]$tr.InsertAfter(", i2c)")
[void]$tr.InsertAfter("`rBrief description of Resources used (")
[void]$tr.InsertAfter("opencv")
[void]$tr.InsertAfter(" and such)")
[void]$tr.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 3 - System Introduction
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "System Introduction"

$tr = $s3.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Introduction of the robot we are building to meet the objective"
[void]$tr.InsertAfter("`rExplain available components and their purpose and the subsystems they are a part of")
[void]$tr.InsertAfter("`rDefine objectives of subsystems")
[void]$tr.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 4 - Motor Subsystem
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Motor Subsystem"

$tr = $s4.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Clearly defined Design process"
[void]$tr.InsertAfter("`rImages and data that show design/build process")
[void]$tr.InsertAfter("`rShow achieved performance of each subsystem")
[void]$tr.InsertAfter("`rGraphs and tables with readable text")
[void]$tr.InsertAfter("`rBlock diagrams described in detail")
[void]$tr.InsertAfter("`rPerformance achievable (goal of subsystem that contributes to overall system success?)")

# ---------------------------------------------------------------------
# Slide 5 - Control Subsystem
# ---------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Control Subsystem"

$tr = $s5.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Clearly defined Design process"
[void]$tr.InsertAfter("`rImages and data that show design/build process")
[void]$tr.InsertAfter("`rShow achieved performance of each subsystem")
[void]$tr.InsertAfter("`rGraphs and tables with readable text")
[void]$tr.InsertAfter("`rBlock diagrams described in detail")
[void]$tr.InsertAfter("`rPerformance achievable (goal of subsystem that contributes to overall system success?)")
[void]$tr.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 6 - Computer Vision Subsystem
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Computer Vision Subsystem"

$tr = $s6.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Clearly defined Design process"
[void]$tr.InsertAfter("`rImages and data that show design/build process")
[void]$tr.InsertAfter("`rShow achieved performance of each subsystem")
[void]$tr.InsertAfter("`rGraphs and tables with readable text")
[void]$tr.InsertAfter("`rBlock diagrams described in detail")
[void]$tr.InsertAfter("`rPerformance achievable (goal of subsystem that contributes to overall system success?)")
[void]$tr.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 7 - Communication and Integration
# ---------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Communication and Integration"

$tr = $s7.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Clearly defined Design process"
[void]$tr.InsertAfter("`rImages and data that show design/build process")
[void]$tr.InsertAfter("`rShow achieved performance of each subsystem")
[void]$tr.InsertAfter("`rGraphs and tables with readable text")
[void]$tr.InsertAfter("`rBlock diagrams described in detail")
[void]$tr.InsertAfter("`rPerformance achievable (goal of subsystem that contributes to overall system success?)")
[void]$tr.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 8 - Simulations and Experiments
# ---------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Simulations and Experiments (each group add something)"

$rsquo = [char]0x2019
$tr = $s8.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Could be part of each subsystem, probably won" + $rsquo + "t be its own piece of the presentation"
